$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A entirely: this shifts columns B:F left to A:E,
# matching the target layout (QS_Phylonet15/FNRATE_EXACT_ASTRAL/TAXON/
# MODEL_CONDITION/GENE header and the associated data), and drops the
# old column A (rank) values/style along with it.
$ws.Range("A:A").Delete()
